$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values that look numeric stay as text, matching source formatting
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.532.15"
$ws.Range("E2").Value = "  -2.88%  "
$ws.Range("D3").Value = "1.654.33"
$ws.Range("E3").Value = "  -4.51%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "214.04"
$ws.Range("D6").Value = "0.509"
$ws.Range("E6").Value = "  -2.48%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").Value = "23.93"
$ws.Range("E8").Value = "  -0.51%  "
$ws.Range("E9").Value = "  -2.63%  "
$ws.Range("D10").Value = "0.0618"
$ws.Range("E10").Value = "  -3.00%  "
$ws.Range("D11").Value = "0.0881"
$ws.Range("E11").Value = "  -1.69%  "
$ws.Range("D12").Value = "1.890.96"
$ws.Range("E12").Value = "  -4.35%  "
$ws.Range("D13").Value = "1.653.41"
$ws.Range("E13").Value = "  -4.50%  "
$ws.Range("D14").Value = "4.13"
$ws.Range("E14").Value = "  -2.93%  "
$ws.Range("D15").Value = "0.562"
$ws.Range("E15").Value = "  -0.31%  "
$ws.Range("D16").Value = "65.77"
$ws.Range("E16").Value = "  -3.03%  "
$ws.Range("D17").Value = "27.522.04"
$ws.Range("E17").Value = "  -2.87%  "
$ws.Range("D18").Value = "240.14"
$ws.Range("E18").Value = "  -2.98%  "
$ws.Range("D19").Value = "0.0₃0728"
$ws.Range("E19").Value = "  -3.57%  "
$ws.Range("D20").Value = "7.53"
$ws.Range("E20").Value = "  -4.84%  "
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").Value = "4.46"
$ws.Range("E22").Value = "  -4.19%  "
$ws.Range("E23").Value = "  -3.96%  "
$ws.Range("E24").Value = "  -2.47%  "
$ws.Range("D25").Value = "145.69"
$ws.Range("E25").Value = "  -2.58%  "
$ws.Range("D26").Value = "7.18"
$ws.Range("E26").Value = "  -4.28%  "
$ws.Range("D27").Value = "16.21"
$ws.Range("E27").Value = "  -3.19%  "
$ws.Range("E28").Value = "  +0.22%  "
$ws.Range("E29").Value = "  -2.56%  "
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").Value = "0.0501"
$ws.Range("E31").Value = "  -2.89%  "
$ws.Range("D32").Value = "3.32"
$ws.Range("E32").Value = "  -3.08%  "
$ws.Range("D33").Value = "1.450.59"
$ws.Range("E33").Value = "  -2.74%  "
$ws.Range("D34").Value = "3.10"
$ws.Range("E34").Value = "  -5.26%  "
$ws.Range("D35").Value = "1.57"
$ws.Range("E35").Value = "  -5.34%  "
$ws.Range("E36").Value = "  -0.81%  "
$ws.Range("D37").Value = "0.919"
$ws.Range("E37").Value = "  -6.16%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.0171"
$ws.Range("E38").Value = "  -3.25%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "0.570"
$ws.Range("E39").Value = "  -5.46%  "
$ws.Range("E40").Value = "  -3.41%  "
$ws.Range("D41").Value = "69.22"
$ws.Range("E41").Value = "  -1.49%  "
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("D43").Value = "5.41"
$ws.Range("E43").Value = "  -4.43%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "0.794"
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("B45").Value = "MXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D45").Value = "2.22"
$ws.Range("E45").Value = "  -3.57%  "
$ws.Range("D46").Value = "1.797.63"
$ws.Range("E46").Value = "  -4.33%  "
$ws.Range("E47").Value = "  -1.40%  "
$ws.Range("D48").Value = "88.40"
$ws.Range("E48").Value = "  -2.29%  "
$ws.Range("D49").Value = "0.0₆0107"
$ws.Range("E49").Value = "  -6.25%  "
$ws.Range("E50").Value = "  -2.18%  "
$ws.Range("D51").Value = "7.80"
$ws.Range("E51").Value = "  -5.18%  "
